$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row for "Short Assignment 3" deadline, mirroring formatting of the
# row above it (e.g. row 14 which holds SA2 / its date, styled as text).
$ws.Range("A15").Value = "SA3"
$ws.Range("B15").Value = "September 23, 2024"

# Match the text number format used by the rest of column B.
$ws.Range("B15").NumberFormat = $ws.Range("B14").NumberFormat

# Update the sheet's active selection to match the saved view state.
$ws.Range("B15").Select()
